$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ("New York -- New York") now has a completed successful API result ---

# B4: Date Published -> 2020-07-22 (serial 44034), keep the same date format used
# by the other rows in column B (yyyy-mm-dd).
$ws.Range("B4").Value = 44034
$ws.Range("B4").NumberFormat = "YYYY-MM-DD"

# C4 / D4: Total Cases / Total Deaths were produced upstream as numeric-looking
# strings (same pattern already seen e.g. in row 8), so force text storage via a
# temporary Text number format, then clear the format again so no stray style
# sticks to the cell.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "219128"
$ws.Range("C4").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "18803"
$ws.Range("D4").ClearFormats()

# E4-H4: numeric counts/percentages
$ws.Range("E4").Value = 33790
$ws.Range("F4").Value = 5239
$ws.Range("G4").Value = 30.07
$ws.Range("H4").Value = 30.43

# J4: Pct Includes Hispanic Black flips to TRUE
$ws.Range("J4").Value = $true

# K4 / L4: Count Cases/Deaths Known Race
$ws.Range("K4").Value = 112360
$ws.Range("L4").Value = 17217

# O4: Status code text - the retried run succeeded this time
$ws.Range("O4").Value = "Success!"

# --- Row 41 ("Iowa"): Count Cases Black/AA corrected ---
$ws.Range("E41").Value = 3288
